$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Bibliography entries for Richtel (2012) and Stross (2010): merge
#    the three consecutive runs that carry the in-text citation /
#    title / closing quote into a single run (same text, same empty
#    formatting) and then apply a red highlight to every run in both
#    paragraphs, including the run inside the hyperlink field.
# ------------------------------------------------------------------

function Merge-CitationRuns($paraIndex, $offset1, $offset2, $offset3, $midText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $run3 = $d.Range($r.Start + $offset2, $r.Start + $offset3)
    $run3.Delete()
    $run2 = $d.Range($r.Start + $offset1, $r.Start + $offset2)
    $run2.Delete()
    $run1 = $d.Range($r.Start + 7, $r.Start + $offset1)
    $run1.InsertAfter($midText)
}

# "Richtel" paragraph - find it by its distinctive text.
$d.Content.Find.ClearFormatting()
$found = $d.Content.Find.Execute("Richtel, M., (2012)")
$richtelPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Richtel, M., (2012)*") {
        $richtelPara = $p
        break
    }
}
$r = $richtelPara.Range
$run3 = $d.Range($r.Start + 62, $r.Start + 65)
$run3.Delete()
$run2 = $d.Range($r.Start + 21, $r.Start + 62)
$run2.Delete()
$run1 = $d.Range($r.Start + 7, $r.Start + 21)
$run1.InsertAfter("Wasting Time Is New Divide in Digital Era" + [char]0x2019 + ", ")

# "Stross" paragraph - find it by its distinctive text.
$strossPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Stross, R. (2010)*") {
        $strossPara = $p
        break
    }
}
$r = $strossPara.Range
$run3 = $d.Range($r.Start + 74, $r.Start + 77)
$run3.Delete()
$run2 = $d.Range($r.Start + 19, $r.Start + 74)
$run2.Delete()
$run1 = $d.Range($r.Start + 6, $r.Start + 19)
$run1.InsertAfter("Computers at Home: Educational Hope vs. Teenage Reality" + [char]0x2019 + ", ")

# Apply a red highlight to the whole of both paragraphs (direct
# character formatting), then patch the hyperlink runs explicitly -
# the hyperlink field's inner run does not pick up highlighting from
# a range that merely spans it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Richtel, M., (2012)*" -or $p.Range.Text -like "Stross, R. (2010)*") {
        $p.Range.Font.HighlightColorIndex = 6
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        foreach ($h in $d.Hyperlinks) {
            if ($h.Range.Start -ge $pStart -and $h.Range.End -le $pEnd) {
                $h.Range.Font.HighlightColorIndex = 6
            }
        }
    }
}

# ------------------------------------------------------------------
# 2. Move the _GoBack bookmark from the "North Carolina study..."
#    paragraph to the very end of the document (end of the Stross
#    bibliography paragraph).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$d.Bookmarks.Add("_GoBack", $endRange)

Write-Host "done"
